$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Fall/Spring/Summer 2025 header block (row 30) and its totals row (row 38),
# along with all the blank rows in between, by deleting the whole row range.
$ws.Range("A30:F38").EntireRow.Delete() | Out-Null

# Fall 2022 / Spring 2022 block: row 7 Spring 2022 course changes from CPSC 6985 (4 credits)
# to CPSC 4115 (3 credits)
$ws.Range("C7").Value = "CPSC 4115"
$ws.Range("D7").Value = 3

# Fall 2023 / Spring 2023 block: Spring 2023 now only has one course (CPSC 6985, 4 credits),
# and the rest of the courses that used to be listed under Spring 2023 move to the
# Fall 2023 column instead.
$ws.Range("C13").Value = "CPSC 6985"
$ws.Range("D13").Value = 4

$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()

$ws.Range("A15").Value = "CYBR 3106"
$ws.Range("B15").Value = 3
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()

$ws.Range("A16").Value = "CYBR 3108"
$ws.Range("B16").Value = 3
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()

$ws.Range("A17").Value = "CYBR 3119"
$ws.Range("B17").Value = 3
